$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "conversion del dia" text block (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$text = $cellA1.Value()
$text = $text.Replace("1000 Bs = 1.62 = 6043.95 pesos", "1000 Bs = 1.62 = 6027.51 pesos")
$text = $text.Replace("6043.95 pesos = 1.62 = 969.16 Bs", "6027.51 pesos = 1.61 = 961.95 Bs")
$cellA1.Value = $text

# --- tasas: update the N/O column rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("O10").Value = 3725
$ws2.Range("N12").Value = 3742
$ws2.Range("O12").Value = 597.2
